$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 0.804
$ws.Range("G3").Value = 0.8080000000000001
$ws.Range("G4").Value = 0.734
$ws.Range("G5").Value = 0.789
$ws.Range("B6").Value = 0.584
$ws.Range("G6").Value = 0.723
$ws.Range("E7").Value = 0.577
$ws.Range("E8").Value = 0.596
$ws.Range("L8").Value = 0.471
$ws.Range("G9").Value = 0.744
$ws.Range("L10").Value = 0.433
$ws.Range("J11").Value = 0.378
$ws.Range("L12").Value = 0.438
$ws.Range("G13").Value = 0.759
$ws.Range("K13").Value = 0.458
$ws.Range("G14").Value = 0.6860000000000001
$ws.Range("G15").Value = 0.756
$ws.Range("J16").Value = 0.337
$ws.Range("H17").Value = 0.361
$ws.Range("J17").Value = 0.325
$ws.Range("H18").Value = 0.36
$ws.Range("I18").Value = 0.5580000000000001
$ws.Range("I19").Value = 0.603
$ws.Range("L19").Value = 0.482
$ws.Range("C20").Value = 0.673
$ws.Range("E20").Value = 0.594
$ws.Range("L21").Value = 0.385
$ws.Range("G22").Value = 0.727
$ws.Range("K22").Value = 0.482
$ws.Range("I23").Value = 0.581
$ws.Range("J23").Value = 0.383
$ws.Range("L24").Value = 0.434
$ws.Range("I25").Value = 0.541
$ws.Range("J25").Value = 0.343
$ws.Range("E26").Value = 0.579
$ws.Range("L27").Value = 0.404
$ws.Range("E28").Value = 0.556
$ws.Range("L29").Value = 0.428
$ws.Range("C30").Value = 0.671
$ws.Range("E30").Value = 0.5590000000000001
$ws.Range("E31").Value = 0.517
$ws.Range("E32").Value = 0.524
$ws.Range("J32").Value = 0.523
$ws.Range("J33").Value = 0.411
$ws.Range("L33").Value = 0.466
$ws.Range("C34").Value = 0.638
$ws.Range("G34").Value = 0.725
$ws.Range("C35").Value = 0.622
$ws.Range("L36").Value = 0.432
$ws.Range("G37").Value = 0.662
$ws.Range("K37").Value = 0.45
$ws.Range("F38").Value = 0.555
$ws.Range("J38").Value = 0.393
$ws.Range("B39").Value = 0.524
$ws.Range("C39").Value = 0.6899999999999999
$ws.Range("L39").Value = 0.389
$ws.Range("B40").Value = 0.571
$ws.Range("I41").Value = 0.548
$ws.Range("J42").Value = 0.361
$ws.Range("D43").Value = 0.494
$ws.Range("H43").Value = 0.428
$ws.Range("E44").Value = 0.573
$ws.Range("H44").Value = 0.387
$ws.Range("C45").Value = 0.662
$ws.Range("G46").Value = 0.638
$ws.Range("L47").Value = 0.453
$ws.Range("D48").Value = 0.473
$ws.Range("J48").Value = 0.411
$ws.Range("E49").Value = 0.533
$ws.Range("D50").Value = 0.496
$ws.Range("E51").Value = 0.551
$ws.Range("G52").Value = 0.661
$ws.Range("L52").Value = 0.432
$ws.Range("H53").Value = 0.411
$ws.Range("J54").Value = 0.383
$ws.Range("C55").Value = 0.6
$ws.Range("H56").Value = 0.434
$ws.Range("H57").Value = 0.367
$ws.Range("L58").Value = 0.431
$ws.Range("E59").Value = 0.575
$ws.Range("E60").Value = 0.517
$ws.Range("D61").Value = 0.456
$ws.Range("L61").Value = 0.433
$ws.Range("I62").Value = 0.612
$ws.Range("J62").Value = 0.48
$ws.Range("E63").Value = 0.542
$ws.Range("F63").Value = 0.494
$ws.Range("G64").Value = 0.772
$ws.Range("E65").Value = 0.603
$ws.Range("L65").Value = 0.447
$ws.Range("D66").Value = 0.538
$ws.Range("D67").Value = 0.515
$ws.Range("D68").Value = 0.497
$ws.Range("G68").Value = 0.729
$ws.Range("F69").Value = 0.554
$ws.Range("G69").Value = 0.671
$ws.Range("I70").Value = 0.635
$ws.Range("E71").Value = 0.554
$ws.Range("B72").Value = 0.634
$ws.Range("G72").Value = 0.658
$ws.Range("G73").Value = 0.675
$ws.Range("L73").Value = 0.46
$ws.Range("F74").Value = 0.47
$ws.Range("L74").Value = 0.405
$ws.Range("E75").Value = 0.523
$ws.Range("I75").Value = 0.501
$ws.Range("I76").Value = 0.614
$ws.Range("J77").Value = 0.474
$ws.Range("C78").Value = 0.644
$ws.Range("H79").Value = 0.41
$ws.Range("J79").Value = 0.414
$ws.Range("J80").Value = 0.361
$ws.Range("J81").Value = 0.511
$ws.Range("L82").Value = 0.41
$ws.Range("E83").Value = 0.556
$ws.Range("G84").Value = 0.697
$ws.Range("H84").Value = 0.412
$ws.Range("I85").Value = 0.631
$ws.Range("E86").Value = 0.546
$ws.Range("J87").Value = 0.404
$ws.Range("G88").Value = 0.72
$ws.Range("J89").Value = 0.415
$ws.Range("B90").Value = 0.641
$ws.Range("G90").Value = 0.667
$ws.Range("E91").Value = 0.51
$ws.Range("C92").Value = 0.669
$ws.Range("H93").Value = 0.444
$ws.Range("J93").Value = 0.476
$ws.Range("H94").Value = 0.401
$ws.Range("E95").Value = 0.516
$ws.Range("L95").Value = 0.395
$ws.Range("G96").Value = 0.72
$ws.Range("J97").Value = 0.512
$ws.Range("B98").Value = 0.676
$ws.Range("J98").Value = 0.433
$ws.Range("B99").Value = 0.724
$ws.Range("H99").Value = 0.454
$ws.Range("D100").Value = 0.544
$ws.Range("H100").Value = 0.474
$ws.Range("H101").Value = 0.345
$ws.Range("L101").Value = 0.416
